# Quarterly indexing esoteric bug-fix operation
#
# The evaluation table rolls forward by one quarter: every existing row's
# metrics (B:G) shift down into the next row, and a brand-new first data
# row (row 2) is populated with the freshly computed error metrics for the
# newest quarter. The oldest row's original values simply fall off the
# bottom of the window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each data row (row 2 = newest quarter results, the rest
# are simply the prior row's old values shifted down by one row).
$data = @{
    2  = @(0.08298509261913574, 0.2782259881569277, 0.1147754943926806, 0.3387853219852958, 0.3408637854649495, 14)
    3  = @(0.1420658796491534,  0.2963740997880928, 0.1213439170623603, 0.3483445378678418, 0.3310457719767695, 13)
    4  = @(0.1568915331285771,  0.2900446020859644, 0.1128329969643262, 0.3359062323987547, 0.3102222847576306, 12)
    5  = @(0.1956053612568337,  0.2558683862338031, 0.08648424928328589,0.2940820451562555, 0.2303151560882531, 11)
    6  = @(0.222511934274951,   0.2817985532879059, 0.1129743571721927, 0.3361165827093223, 0.2655451338021729, 10)
    7  = @(0.176450614684246,   0.1849112350951672, 0.06089664821255769,0.2467724624275523, 0.1829810301339271, 9)
    8  = @(0.2112811930113801,  0.285830133390834,  0.1115765264066959, 0.3340307267403643, 0.2765848542209762, 8)
    9  = @(0.1601567932497119,  0.2716627869588349, 0.08940684574059213,0.2990097753261457, 0.2727320208605046, 7)
    10 = @(0.2627214152773011,  0.3235434349789264, 0.1669924024756296, 0.4086470389904099, 0.3380801243818421, 7)
    11 = @(0.07165448440387559, 0.1716308632352222, 0.0355562541576917, 0.1885636607559678, 0.1950060544653244, 5)
}

$columns = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
